$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 3 data rows (old rows 2-4), shifting rows 5-21 up to 2-18
$ws.Rows("2:4").Delete()

# Append 13 new data rows (rows 19-31) with the May 9th measurements
$ws.Range("A19").Value = 5.261263814465728
$ws.Range("B19").Value = -3.646456844505198
$ws.Range("C19").Value = -4.444081427037008
$ws.Range("A20").Value = -0.0522782514835689
$ws.Range("B20").Value = 0.8027088176245329
$ws.Range("C20").Value = -0.5144340389076021
$ws.Range("A21").Value = -7.4421002289345
$ws.Range("B21").Value = -0.806786348079851
$ws.Range("C21").Value = 5.891316611191369
$ws.Range("A22").Value = -0.6881247882184418
$ws.Range("B22").Value = -10.85945387133243
$ws.Range("C22").Value = 4.675690663957009
$ws.Range("A23").Value = 4.403560649389491
$ws.Range("B23").Value = -11.45763061786516
$ws.Range("C23").Value = 1.461910155997879
$ws.Range("A24").Value = 4.784720346845424
$ws.Range("B24").Value = 3.472261708358207
$ws.Range("C24").Value = -1.690446103441295
$ws.Range("A25").Value = 0.8494467159797106
$ws.Range("B25").Value = 5.416543818924583
$ws.Range("C25").Value = 0.3924825684777868
$ws.Range("A26").Value = -4.816584700825577
$ws.Range("B26").Value = 1.840564275695681
$ws.Range("C26").Value = 3.20516648785821
$ws.Range("A27").Value = -4.067581341184413
$ws.Range("B27").Value = 7.491072893142713
$ws.Range("C27").Value = -0.0797111165934683
$ws.Range("A28").Value = -0.603338341945886
$ws.Range("B28").Value = 10.07777972605036
$ws.Range("C28").Value = -4.32223105156559
$ws.Range("A29").Value = 3.495839129919295
$ws.Range("B29").Value = 0.1920398131184147
$ws.Range("C29").Value = -4.956773451004905
$ws.Range("A30").Value = 4.441888874974746
$ws.Range("B30").Value = -4.70699503503997
$ws.Range("C30").Value = -3.40560439948378
$ws.Range("A31").Value = 1.070874737597066
$ws.Range("B31").Value = -0.5835215826144202
$ws.Range("C31").Value = -1.073033490400209
